$wb = $excel.ActiveWorkbook

# The Overview sheet mirrors the per-language status in row 2 (the md file row)
# via the same shared string, so it needs to reflect the new status text too.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handoff transform failed"
$overview.Range("C2").Value = "Handoff transform failed"

$sheetNames = @("zh-cn", "de-de")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Remove the hyperlink attached to C2 (the "Latest Handoff File" cell) only,
    # leaving the other hyperlinks (A2, A3) intact.
    $targetLink = $null
    foreach ($hl in $ws.Hyperlinks) {
        $cell = $hl.Range.Cells.Item(1, 1)
        if ($cell.Row -eq 2 -and $cell.Column -eq 3) {
            $targetLink = $hl
        }
    }
    if ($targetLink -ne $null) {
        $targetLink.Delete()
    }

    # Fully clear the now-unlinked "Latest Handoff File" cell (value + formatting).
    $ws.Range("C2").Clear()

    # Update status / reason / handoff datetime to reflect a failed handoff.
    $ws.Range("B2").Value = "Handoff transform failed"
    $ws.Range("D2").Value = "0001-01-01 00:00:00"
    $ws.Range("H2").Value = "Ignored"
}
